$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRange, $value) {
    $cellRange.NumberFormat = "@"
    $cellRange.Value = $value
    $cellRange.ClearFormats()
}

Set-TextValue $ws.Range("D2") "56.563.68"
Set-TextValue $ws.Range("E2") "  -3.93%  "
Set-TextValue $ws.Range("D3") "2.358.89"
Set-TextValue $ws.Range("E3") "  -5.37%  "
Set-TextValue $ws.Range("E4") "  -0.33%  "
Set-TextValue $ws.Range("D5") "513.50"
Set-TextValue $ws.Range("D6") "127.80"
Set-TextValue $ws.Range("E6") "  -5.84%  "
Set-TextValue $ws.Range("D7") "0.998"
Set-TextValue $ws.Range("E7") "  -0.31%  "
Set-TextValue $ws.Range("E8") "  -2.26%  "
Set-TextValue $ws.Range("D9") "2.374.13"
Set-TextValue $ws.Range("E9") "  -5.73%  "
Set-TextValue $ws.Range("D10") "0.0960"
Set-TextValue $ws.Range("E10") "  -3.73%  "
Set-TextValue $ws.Range("E11") "  -1.72%  "
Set-TextValue $ws.Range("E12") "  -8.60%  "
Set-TextValue $ws.Range("E13") "  -5.76%  "
Set-TextValue $ws.Range("D14") "2.777.88"
Set-TextValue $ws.Range("E14") "  -6.34%  "
Set-TextValue $ws.Range("D15") "56.480.08"
Set-TextValue $ws.Range("E15") "  -4.03%  "
Set-TextValue $ws.Range("D16") "21.44"
Set-TextValue $ws.Range("E16") "  -4.83%  "
Set-TextValue $ws.Range("E17") "  -4.37%  "
Set-TextValue $ws.Range("D18") "2.352.35"
Set-TextValue $ws.Range("E18") "  -6.59%  "
Set-TextValue $ws.Range("E19") "  -4.15%  "
Set-TextValue $ws.Range("D20") "309.53"
Set-TextValue $ws.Range("E20") "  -3.96%  "
Set-TextValue $ws.Range("E21") "  -5.46%  "
Set-TextValue $ws.Range("D22") "6.09"
Set-TextValue $ws.Range("E22") "  -1.01%  "
Set-TextValue $ws.Range("D23") "0.997"
Set-TextValue $ws.Range("E23") "  -0.03%  "
Set-TextValue $ws.Range("D24") "64.67"
Set-TextValue $ws.Range("D25") "1.00"
Set-TextValue $ws.Range("E25") "  +0.39%  "
Set-TextValue $ws.Range("E26") "  -4.94%  "
Set-TextValue $ws.Range("D27") "2.466.79"
Set-TextValue $ws.Range("E27") "  -6.60%  "
Set-TextValue $ws.Range("E28") "  -4.56%  "
Set-TextValue $ws.Range("D29") "7.18"
Set-TextValue $ws.Range("E29") "  -4.89%  "
Set-TextValue $ws.Range("D30") "174.04"
Set-TextValue $ws.Range("E30") "  +1.21%  "
Set-TextValue $ws.Range("E31") "  -5.25%  "
Set-TextValue $ws.Range("D32") "0.0₃0715"
Set-TextValue $ws.Range("E32") "  -6.75%  "
Set-TextValue $ws.Range("D33") "6.10"
Set-TextValue $ws.Range("E33") "  -4.18%  "
Set-TextValue $ws.Range("D34") "1.13"
Set-TextValue $ws.Range("E34") "  -7.22%  "
Set-TextValue $ws.Range("E35") "  -0.06%  "
Set-TextValue $ws.Range("E36") "  +0.05%  "
Set-TextValue $ws.Range("E37") "  -3.58%  "
Set-TextValue $ws.Range("E38") "  -6.00%  "
Set-TextValue $ws.Range("D39") "3.72"
Set-TextValue $ws.Range("E39") "  -7.14%  "
Set-TextValue $ws.Range("E40") "  +1.50%  "
Set-TextValue $ws.Range("E42") "  -6.81%  "
Set-TextValue $ws.Range("E43") "  -4.92%  "
Set-TextValue $ws.Range("D44") "4.87"
Set-TextValue $ws.Range("E44") "  -4.58%  "
Set-TextValue $ws.Range("E45") "  -5.04%  "
Set-TextValue $ws.Range("D46") "122.01"
Set-TextValue $ws.Range("E46") "  -7.60%  "
Set-TextValue $ws.Range("D47") "252.12"
Set-TextValue $ws.Range("E47") "  -10.26%  "
Set-TextValue $ws.Range("D48") "0.0904"
Set-TextValue $ws.Range("E48") "  -2.85%  "
Set-TextValue $ws.Range("E49") "  -4.69%  "
Set-TextValue $ws.Range("E50") "  -5.81%  "
Set-TextValue $ws.Range("D51") "16.61"
Set-TextValue $ws.Range("E51") "  -6.50%  "
